# "Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta"
#
# Updates the "Estado de Cuenta" (account statement) table on sheet Hoja1:
#  - Row 16/17 (GREGORIO ROBLEDO CASTRO, periods 1609 / 2508): the
#    "Periodo Mora" labels and their matching "Valor Mora" amounts are
#    corrected/swapped to the right period.
#  - Rows 18/19 (YEIMY YOHANA BOTERO MONTOYA / SANTIAGO URREGO BOTERO):
#    the "Periodo Mora" period is updated from 2507 to 2508.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Periodo Mora (column E)
$ws.Range("E16").Value = "1609"
$ws.Range("E17").Value = "2508"
$ws.Range("E18").Value = "2508"
$ws.Range("E19").Value = "2508"

# Valor Mora (column F) - values for rows 16 and 17 swap along with the
# corrected periods above
$ws.Range("F16").Value = 27578
$ws.Range("F17").Value = 227760
